$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "MONTO MMOO" column header in G1, matching the centered
# header style already used by A1:F1.
$ws.Range("G1").Value = "MONTO MMOO"
$ws.Range("G1").HorizontalAlignment = -4108

# Give the new column a sensible width, matching the other bestFit columns.
$ws.Columns.Item(7).ColumnWidth = 14.7109375

# The old blank E/F/G cells (rows 2-10) only ever carried an unused
# "applyBorder" style - remove them so the cells go back to the sheet's
# implicit default formatting.
$ws.Range("E2:G10").Clear()

# Mark the new MMOO amount cell (G6) with an underlined font, same spot
# referenced by the updated selection below.
$ws.Range("G6").Font.Underline = $true

# Drop the three trailing all-default rows that no longer belong to the
# table body.
$ws.Rows("11:13").Delete()

# Match the new active selection left behind by the edit.
$ws.Range("G6").Select()
